# "updated main GSC export data"
# Append the next day's row (2025-12-02) to the bottom of the "Chart" sheet,
# continuing the existing Date / Invalid / Valid series with the same
# values as the prior day (Invalid = 0, Valid = 27).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$lastRow = $ws.UsedRange.Rows.Count
$newRow = $lastRow + 1

# Force the new date to be stored as literal text (matching the existing
# "Date" column, which holds text like "2025-12-01" rather than real date
# serials), then clear the number-format bit so the cell keeps the
# worksheet's default (General) style instead of picking up a date format.
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "2025-12-02"
$ws.Cells.Item($newRow, 1).ClearFormats()

$ws.Cells.Item($newRow, 2).Value = 0
$ws.Cells.Item($newRow, 3).Value = 27
